$d = $word.ActiveDocument

# Fix the "Annexe 1" mistake: the paragraph about the class relational schema
# should reference Annexe 2 (not Annexe 1, which is already used by the
# database relational schema paragraph).
$d.Content.Find.Execute(
    "des classes qui interagissent avec la base de données a été créé (voir annexe 1).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "des classes qui interagissent avec la base de données a été créé (voir annexe 2).",
    2
)
